$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 699.8
$ws.Range("I11").Value = 699.8
$ws.Range("K11").Value = 699.8
$ws.Range("M11").Value = -559.8
$ws.Range("H41").Value = 452.7143
$ws.Range("I41").Value = 436.41666
$ws.Range("J41").Value = 550.5
$ws.Range("K41").Value = 436.41666
$ws.Range("L41").Value = 550.5
$ws.Range("M41").Value = 3.583340000000021
$ws.Range("N41").Value = -1430.5
$ws.Range("H113").Value = 998.8
$ws.Range("J113").Value = 999.75
$ws.Range("L113").Value = 999.75
$ws.Range("N113").Value = -7507.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H32").Value = 1406
$ws.Range("I32").Value = 1406
$ws.Range("K32").Value = 1406
$ws.Range("M32").Value = -1119
$ws.Range("H61").Value = 2757.125
$ws.Range("I61").Value = 2828.6
$ws.Range("J61").Value = 2638
$ws.Range("K61").Value = 2828.6
$ws.Range("L61").Value = 2638
$ws.Range("M61").Value = -2616.6
$ws.Range("N61").Value = -3062
$ws.Range("H74").Value = 1991.1666
$ws.Range("I74").Value = 1986.75
$ws.Range("K74").Value = 1986.75
$ws.Range("M74").Value = -1112.75
$ws.Range("H77").Value = 1991.1666
$ws.Range("I77").Value = 1986.75
$ws.Range("K77").Value = 9933.75
$ws.Range("M77").Value = -5565.75
$ws.Range("H124").Value = 47467
$ws.Range("J124").Value = 47467
$ws.Range("L124").Value = 47467
$ws.Range("N124").Value = -57287
$ws.Range("H136").Value = 2757.125
$ws.Range("I136").Value = 2828.6
$ws.Range("J136").Value = 2638
$ws.Range("K136").Value = 8485.799999999999
$ws.Range("L136").Value = 7914
$ws.Range("M136").Value = -5935.799999999999
$ws.Range("N136").Value = -13014
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 696.8333
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -846
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 737.1667
$ws.Range("I22").Value = 824.375
$ws.Range("K22").Value = 824.375
$ws.Range("M22").Value = -474.375
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("H31").Value = 4641.5
$ws.Range("I31").Value = 3784.4
$ws.Range("J31").Value = 5498.6
$ws.Range("K31").Value = 3784.4
$ws.Range("L31").Value = 5498.6
$ws.Range("M31").Value = -3489.4
$ws.Range("N31").Value = -6088.6
$ws.Range("H34").Value = 4641.5
$ws.Range("I34").Value = 3784.4
$ws.Range("J34").Value = 5498.6
$ws.Range("K34").Value = 3784.4
$ws.Range("L34").Value = 5498.6
$ws.Range("M34").Value = -3582.4
$ws.Range("N34").Value = -5902.6
$ws.Range("H132").Value = 846.75
$ws.Range("I132").Value = 846.75
$ws.Range("K132").Value = 2540.25
$ws.Range("M132").Value = -10.25
$ws.Range("H134").Value = 2159.111
$ws.Range("I134").Value = 2061.7144
$ws.Range("K134").Value = 6185.1432
$ws.Range("M134").Value = -3650.1432
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 1030
$ws.Range("I28").Value = 1030
$ws.Range("K28").Value = 3090
$ws.Range("M28").Value = -2858
$ws.Range("H45").Value = 4000
$ws.Range("J45").Value = 4000
$ws.Range("L45").Value = 12000
$ws.Range("N45").Value = -13064
$ws.Range("H108").Value = 400.57144
$ws.Range("I108").Value = 400.57144
$ws.Range("K108").Value = 1201.71432
$ws.Range("M108").Value = 1678.28568
$ws.Range("H109").Value = 1612.5
$ws.Range("I109").Value = 483.33334
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 1450.00002
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -410.0000199999999
$ws.Range("N109").Value = -17080
$ws.Range("H117").Value = 1283.7778
$ws.Range("I117").Value = 327
$ws.Range("J117").Value = 2049.2
$ws.Range("K117").Value = 981
$ws.Range("L117").Value = 6147.599999999999
$ws.Range("M117").Value = 2461
$ws.Range("N117").Value = -13031.6
$ws.Range("H131").Value = 1923.1818
$ws.Range("J131").Value = 3498.75
$ws.Range("L131").Value = 10496.25
$ws.Range("N131").Value = -20576.25
$ws.Range("H137").Value = 2033
$ws.Range("J137").Value = 2033
$ws.Range("L137").Value = 6099
$ws.Range("N137").Value = -16299
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 334.64706
$ws.Range("I2").Value = 350.6
$ws.Range("J2").Value = 215
$ws.Range("K2").Value = 350.6
$ws.Range("L2").Value = 215
$ws.Range("M2").Value = -237.6
$ws.Range("N2").Value = -441
$ws.Range("H132").Value = 2642
$ws.Range("I132").Value = 1956
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 5868
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -3338
$ws.Range("N132").Value = -17102
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""
$ws.Range("H16").Value = 1179.7778
$ws.Range("I16").Value = 1179.7778
$ws.Range("K16").Value = 1179.7778
$ws.Range("M16").Value = -1009.7778
$ws.Range("H35").Value = 10096.5
$ws.Range("I35").Value = 1495
$ws.Range("K35").Value = 1495
$ws.Range("M35").Value = -1159
$ws.Range("H40").Value = 3476
$ws.Range("I40").Value = 3476
$ws.Range("K40").Value = 3476
$ws.Range("M40").Value = -3340
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = ""
$ws.Range("H136").Value = 3533.3333
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -13500
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 45000
$ws.Range("I41").Value = 45000
$ws.Range("K41").Value = 45000
$ws.Range("M41").Value = -44610
$ws.Range("H108").Value = 67499.5
$ws.Range("J108").Value = 67499.5
$ws.Range("L108").Value = 67499.5
$ws.Range("N108").Value = -75179.5
$ws.Range("H136").Value = 1556.875
$ws.Range("I136").Value = 1064.2858
$ws.Range("K136").Value = 3192.8574
$ws.Range("M136").Value = -642.8574000000003
